$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 83, shifting existing rows 83-99 down to 84-100
$ws.Rows("83:83").Insert()

# Populate the newly inserted row 83 with the new data record
$ws.Range("A83").Value = 3
$ws.Range("B83").Value = "Femacal de La Calera"
$ws.Range("C83").Value = "Coquimbo"
$ws.Range("D83").Value = 45218
$ws.Range("E83").Value = 5
$ws.Range("F83").Value = 100112022
$ws.Range("G83").Value = "Arveja Verde"
$ws.Range("H83").Value = "Perfection"
$ws.Range("I83").Value = "Primera"
$ws.Range("J83").Value = 36
$ws.Range("K83").Value = 27000
$ws.Range("L83").Value = 27000
$ws.Range("M83").Value = 27000
$ws.Range("N83").Value = "$/saco 25 kilos"
$ws.Range("O83").Value = "Provincia de Limarí"
$ws.Range("P83").Value = 1080
$ws.Range("Q83").Value = 25
$ws.Range("R83").Value = "Hortaliza"
